$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at BQ; source_material_id (BP) stays, everything from the
# old BQ (store_cond) onward shifts one column to the right.
$ws.Columns("BQ:BQ").Insert()

# Give the new column its header text (new shared string "specimen_voucher").
$ws.Range("BQ15").Value = "specimen_voucher"

# Column-insert does not relocate cell comments (notes) in this engine, so the
# existing notes are still anchored to their old letters (BQ15.."explain how...",
# BR15.."Information about the genetic distinctness...", etc). Walk the range from
# the last column backwards, moving each note one column to the right so the note
# content tracks the data that moved with it.
$t = 'reference or method used in determining the water content of soil'
$ws.Range("CB15").Comment.Delete()
$ws.Range("CC15").AddComment($t)
$t = 'water content (g/g or cm3/cm3)'
$ws.Range("CA15").Comment.Delete()
$ws.Range("CB15").AddComment($t)
$t = 'Feeding position in food chain (eg., chemolithotroph)'
$ws.Range("BZ15").Comment.Delete()
$ws.Range("CA15").AddComment($t)
$t = 'Definition for soil: total organic C content of the soil units of g C/kg soil. Definition otherwise: total organic carbon content'
$ws.Range("BY15").Comment.Delete()
$ws.Range("BZ15").AddComment($t)
$t = 'reference or method used in determining total organic C'
$ws.Range("BX15").Comment.Delete()
$ws.Range("BY15").AddComment($t)
$t = 'total nitrogen content of the sample'
$ws.Range("BW15").Comment.Delete()
$ws.Range("BX15").AddComment($t)
$t = 'reference or method used in determining the total N'
$ws.Range("BV15").Comment.Delete()
$ws.Range("BW15").AddComment($t)
$t = 'note method(s) used for tilling'
$ws.Range("BU15").Comment.Delete()
$ws.Range("BV15").AddComment($t)
$t = 'reference or method used in determining soil texture'
$ws.Range("BT15").Comment.Delete()
$ws.Range("BU15").AddComment($t)
$t = 'the relative proportion of different grain sizes of mineral particles in a soil, as described using a standard system; express as % sand (50 um to 2 mm), silt (2 um to 50 um), and clay (<2 um) with textural name (e.g., silty clay loam) optional.'
$ws.Range("BS15").Comment.Delete()
$ws.Range("BT15").AddComment($t)
$t = 'Information about the genetic distinctness of the lineage (eg., biovar, serovar)'
$ws.Range("BR15").Comment.Delete()
$ws.Range("BS15").AddComment($t)
$t = 'explain how and for how long the soil sample was stored before DNA extraction.'
$ws.Range("BQ15").Comment.Delete()
$ws.Range("BR15").AddComment($t)

# Finally, give the (now vacated) BQ15 header its own note describing specimen_voucher.
$specimenNote = 'Identifier for the physical specimen. Use format: "[<institution-code>:[<collection-code>:]]<specimen_id>", eg, "UAM:Mamm:52179". Intended as a reference to the physical specimen that remains after it was analyzed. If the specimen was destroyed in the process of analysis, electronic images (e-vouchers) are an adequate substitute for a physical voucher specimen. Ideally the specimens will be deposited in a curated museum, herbarium, or frozen tissue collection, but often they will remain in a personal or laboratory collection for some time before they are deposited in a curated collection. There are three forms of specimen_voucher qualifiers. If the text of the qualifier includes one or more colons it is a ''structured voucher''. Structured vouchers include institution-codes (and optional collection-codes) taken from a controlled vocabulary maintained by the INSDC that denotes the museum or herbarium collection where the specimen resides, please visit the INSDC website, http://www.insdc.org/controlled-vocabulary-specimenvoucher-qualifier'
$ws.Range("BQ15").AddComment($specimenNote)

Write-Host "done"
